$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 240.66667
$ws.Range("J33").Value = 173.66667
$ws.Range("L33").Value = 173.66667
$ws.Range("N33").Value = -631.6666700000001
$ws.Range("H38").Value = 1909.1333
$ws.Range("I38").Value = 615.6
$ws.Range("J38").Value = 4496.2
$ws.Range("K38").Value = 1846.8
$ws.Range("L38").Value = 13488.6
$ws.Range("M38").Value = -1474.8
$ws.Range("N38").Value = -14232.6
$ws.Range("H74").Value = 6569.154
$ws.Range("I74").Value = 6569.154
$ws.Range("K74").Value = 6569.154
$ws.Range("M74").Value = -5633.154
$ws.Range("H77").Value = 6569.154
$ws.Range("I77").Value = 6569.154
$ws.Range("K77").Value = 32845.77
$ws.Range("M77").Value = -28165.77
$ws.Range("H86").Value = 3437.2632
$ws.Range("I86").Value = 3246.889
$ws.Range("K86").Value = 3246.889
$ws.Range("M86").Value = -2123.889
$ws.Range("H89").Value = 3437.2632
$ws.Range("I89").Value = 3246.889
$ws.Range("K89").Value = 16234.445
$ws.Range("M89").Value = -10618.445
$ws.Range("H107").Value = 3340.8462
$ws.Range("I107").Value = 2040.2727
$ws.Range("K107").Value = 2040.2727
$ws.Range("M107").Value = -120.2727
$ws.Range("H113").Value = 2661.25
$ws.Range("I113").Value = 1951.25
$ws.Range("J113").Value = 3371.25
$ws.Range("K113").Value = 1951.25
$ws.Range("L113").Value = 3371.25
$ws.Range("M113").Value = 1302.75
$ws.Range("N113").Value = -9879.25
$ws.Range("H137").Value = 2929.4092
$ws.Range("I137").Value = 2737.423
$ws.Range("K137").Value = 8212.269
$ws.Range("M137").Value = -5662.269

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4402.2856
$ws.Range("I32").Value = 3427.6128
$ws.Range("K32").Value = 3427.6128
$ws.Range("M32").Value = -3140.6128
$ws.Range("H48").Value = 126717
$ws.Range("J48").Value = 126717
$ws.Range("L48").Value = 126717
$ws.Range("N48").Value = -127485
$ws.Range("H61").Value = 1928.4286
$ws.Range("I61").Value = 1874.75
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1874.75
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1662.75
$ws.Range("N61").Value = -2424
$ws.Range("H122").Value = 6622
$ws.Range("I122").Value = 5035.048
$ws.Range("J122").Value = 7119.403
$ws.Range("K122").Value = 15105.144
$ws.Range("L122").Value = 21358.209
$ws.Range("M122").Value = -12655.144
$ws.Range("N122").Value = -26258.209
$ws.Range("H132").Value = 5260.8076
$ws.Range("I132").Value = 2620.2
$ws.Range("K132").Value = 7860.599999999999
$ws.Range("M132").Value = -5330.599999999999
$ws.Range("H136").Value = 1928.4286
$ws.Range("I136").Value = 1874.75
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5624.25
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3074.25
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 10000
$ws.Range("I7").Value = 10000
$ws.Range("K7").Value = 10000
$ws.Range("M7").Value = -9887
$ws.Range("H20").Value = 2793.9534
$ws.Range("I20").Value = 2607.05
$ws.Range("K20").Value = 2607.05
$ws.Range("M20").Value = -2360.05
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H43").Value = 99999
$ws.Range("J43").Value = 99999
$ws.Range("L43").Value = 99999
$ws.Range("N43").Value = -100361
$ws.Range("H47").Value = 87120.336
$ws.Range("J47").Value = 87120.336
$ws.Range("L47").Value = 87120.336
$ws.Range("N47").Value = -88160.336
$ws.Range("H48").Value = 93776.664
$ws.Range("J48").Value = 93776.664
$ws.Range("L48").Value = 93776.664
$ws.Range("N48").Value = -94606.664
$ws.Range("H107").Value = 3480.5715
$ws.Range("I107").Value = 2592.6843
$ws.Range("J107").Value = 5355
$ws.Range("K107").Value = 2592.6843
$ws.Range("L107").Value = 5355
$ws.Range("M107").Value = -672.6842999999999
$ws.Range("N107").Value = -9195

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 125003320
$ws.Range("I62").Value = 166669500
$ws.Range("K62").Value = 166669500
$ws.Range("M62").Value = -166668876
$ws.Range("H65").Value = 125003320
$ws.Range("I65").Value = 166669500
$ws.Range("K65").Value = 833347500
$ws.Range("M65").Value = -833344380
$ws.Range("H94").Value = 3312.75
$ws.Range("J94").Value = 3612.5
$ws.Range("L94").Value = 3612.5
$ws.Range("N94").Value = -4514.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 52.92
$ws.Range("I12").Value = 50.8
$ws.Range("J12").Value = 54.333332
$ws.Range("K12").Value = 152.4
$ws.Range("L12").Value = 162.999996
$ws.Range("M12").Value = 20.60000000000002
$ws.Range("N12").Value = -508.999996
$ws.Range("H14").Value = 11229.134
$ws.Range("I14").Value = 11229.134
$ws.Range("K14").Value = 33687.402
$ws.Range("M14").Value = -33514.402
$ws.Range("H92").Value = 479.96875
$ws.Range("J92").Value = 511.25
$ws.Range("L92").Value = 1533.75
$ws.Range("N92").Value = -4029.75
$ws.Range("H111").Value = 1950.8
$ws.Range("I111").Value = 627
$ws.Range("J111").Value = 2833.3333
$ws.Range("K111").Value = 1881
$ws.Range("L111").Value = 8499.999899999999
$ws.Range("M111").Value = 1186
$ws.Range("N111").Value = -14633.9999
$ws.Range("H119").Value = 4590
$ws.Range("I119").Value = 1180.6666
$ws.Range("K119").Value = 3541.9998
$ws.Range("M119").Value = 1296.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 45000
$ws.Range("J74").Value = 45000
$ws.Range("L74").Value = 45000
$ws.Range("N74").Value = -46872
$ws.Range("H75").Value = 45000
$ws.Range("J75").Value = 45000
$ws.Range("L75").Value = 45000
$ws.Range("N75").Value = -46748
$ws.Range("H77").Value = 45000
$ws.Range("J77").Value = 45000
$ws.Range("L77").Value = 135000
$ws.Range("N77").Value = -144360
$ws.Range("H78").Value = 45000
$ws.Range("J78").Value = 45000
$ws.Range("L78").Value = 135000
$ws.Range("N78").Value = -143736
$ws.Range("H94").Value = 43999
$ws.Range("J94").Value = 43999
$ws.Range("L94").Value = 43999
$ws.Range("N94").Value = -45351
$ws.Range("H132").Value = 8330.895
$ws.Range("I132").Value = 7344.364
$ws.Range("K132").Value = 22033.092
$ws.Range("M132").Value = -19503.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 11610.577
$ws.Range("I61").Value = 11970.208
$ws.Range("K61").Value = 11970.208
$ws.Range("M61").Value = -11768.208
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H93").Value = 3944.4915
$ws.Range("I93").Value = 4203.1953
$ws.Range("J93").Value = 3355.2222
$ws.Range("K93").Value = 4203.1953
$ws.Range("L93").Value = 3355.2222
$ws.Range("M93").Value = -2955.1953
$ws.Range("N93").Value = -5851.2222
$ws.Range("H113").Value = 11610.577
$ws.Range("I113").Value = 11970.208
$ws.Range("K113").Value = 11970.208
$ws.Range("M113").Value = -9800.208000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 54356.07
$ws.Range("J46").Value = 54356.07
$ws.Range("L46").Value = 54356.07
$ws.Range("N46").Value = -54818.07
$ws.Range("H96").Value = 6362.5
$ws.Range("I96").Value = 6362.5
$ws.Range("K96").Value = 6362.5
$ws.Range("M96").Value = -4989.5
$ws.Range("H100").Value = 1012.0645
$ws.Range("I100").Value = 523.4737
$ws.Range("K100").Value = 1046.9474
$ws.Range("M100").Value = -505.9474
$ws.Range("H107").Value = 1505.5312
$ws.Range("I107").Value = 1545.7
$ws.Range("K107").Value = 4637.1
$ws.Range("M107").Value = -2717.1
$ws.Range("H134").Value = 54356.07
$ws.Range("J134").Value = 54356.07
$ws.Range("L134").Value = 163068.21
$ws.Range("N134").Value = -168138.21
$ws.Range("H136").Value = 6851.5654
$ws.Range("I136").Value = 10582.5
$ws.Range("J136").Value = 2781.4546
$ws.Range("K136").Value = 31747.5
$ws.Range("L136").Value = 8344.363799999999
$ws.Range("M136").Value = -29197.5
$ws.Range("N136").Value = -13444.3638
